# Daily attendance processing - 2026-01-01 06:44:07
# Normalize the "Recorded By" (column G) cells: when a cell lists multiple
# recorders separated by ", ", the last two recorders get swapped - except
# for the specific "admin@admin.com, System" pairing, which is left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value2

    if ([string]::IsNullOrEmpty($val)) {
        continue
    }

    if ($val -eq "admin@admin.com, System") {
        continue
    }

    $parts = $val -split ", "
    if ($parts.Count -ge 2) {
        $last = $parts[$parts.Count - 1]
        $secondLast = $parts[$parts.Count - 2]
        $parts[$parts.Count - 1] = $secondLast
        $parts[$parts.Count - 2] = $last
        $cell.Value = [string]::Join(", ", $parts)
    }
}
